$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 16, shifting the "ChatVRM" entry (and
# everything below it) down by one row. This matches the diff: row 16
# becomes row 17, and the sheet dimension/selection grows to A1:A17.
$ws.Rows("16").Insert()

$ws.Range("A17").Select()
